$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Accidente"
$ws.Range("C3").Value = "damn"
$ws.Range("D3").Value = 4.858047097862472
$ws.Range("E3").Value = -74.07229806514576
$ws.Range("F3").Value = "2025-11-13 14:17:08"
$ws.Range("G3").Value = 0
